$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.997.26"
$ws.Range("E2").Value = "  +0.94%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.299.13"
$ws.Range("E3").Value = "  +0.50%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "507.92"
$ws.Range("E5").Value = "  +0.80%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.17"
$ws.Range("E6").Value = "  -0.12%  "

# Row 7
$ws.Range("E7").Value = "  -0.32%  "

# Row 8
$ws.Range("E8").Value = "  +0.31%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.324.94"
$ws.Range("E9").Value = "  +1.13%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0981"
$ws.Range("E10").Value = "  +2.22%  "

# Row 11
$ws.Range("E11").Value = "  +1.68%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.06"
$ws.Range("E12").Value = "  +7.15%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.02"
$ws.Range("E14").Value = "  +4.34%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.709.97"
$ws.Range("E15").Value = "  +0.49%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "54.924.68"
$ws.Range("E16").Value = "  +0.89%  "

# Row 17
$ws.Range("E17").Value = "  +1.47%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.290.70"
$ws.Range("E18").Value = "  -0.32%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.71"
$ws.Range("E19").Value = "  +3.81%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.19"
$ws.Range("E20").Value = "  +0.62%  "

# Row 21
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.69"
$ws.Range("E21").Value = "  +4.19%  "

# Row 22
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "311.80"
$ws.Range("E22").Value = "  +2.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.40%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.53"
$ws.Range("E24").Value = "  -2.33%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.992"
$ws.Range("E25").Value = "  -0.53%  "

# Row 26
$ws.Range("E26").Value = "  -0.15%  "

# Row 27
$ws.Range("E27").Value = "  +2.36%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "173.02"
$ws.Range("E28").Value = "  -0.13%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.17"
$ws.Range("E29").Value = "  +2.71%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0712"
$ws.Range("E30").Value = "  +2.41%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.15"
$ws.Range("E31").Value = "  +4.58%  "

# Row 32
$ws.Range("E32").Value = "  +0.29%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.11"
$ws.Range("E33").Value = "  +1.30%  "

# Row 34
$ws.Range("E34").Value = "  -0.03%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.993"
$ws.Range("E35").Value = "  -0.23%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.23"
$ws.Range("E36").Value = "  +2.11%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.918"
$ws.Range("E37").Value = "  -5.56%  "

# Row 38
$ws.Range("E38").Value = "  +3.24%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.80"
$ws.Range("E39").Value = "  +2.02%  "

# Row 40
$ws.Range("E40").Value = "  +1.65%  "

# Row 41
$ws.Range("E41").Value = "  +0.81%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "135.96"
$ws.Range("E42").Value = "  +8.01%  "

# Row 43
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.44"
$ws.Range("E43").Value = "  +1.00%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.93"
$ws.Range("E44").Value = "  +0.71%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "261.84"
$ws.Range("E45").Value = "  +6.89%  "

# Row 46
$ws.Range("E46").Value = "  +1.48%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0914"
$ws.Range("E47").Value = "  +1.98%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.556"
$ws.Range("E48").Value = "  +1.16%  "

# Row 49
$ws.Range("E49").Value = "  +1.16%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0210"
$ws.Range("E50").Value = "  +1.60%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.83"
$ws.Range("E51").Value = "  +0.42%  "
